# Insert a new data row before the current row 495 ("Hortaliza, Feria Lagunitas
# de Puerto Montt - Betarraga" sheet). This pushes the existing rows 495-544
# down to 496-545 and extends the used range to A1:R545.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(495).Insert()

# Fill in the newly inserted row with the new weekly price record.
$ws.Range("A495").Value = 4
$ws.Range("B495").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C495").Value = "Los Lagos"
$ws.Range("D495").Value = 45212
$ws.Range("E495").Value = 10
$ws.Range("F495").Value = 100114014
$ws.Range("G495").Value = "Betarraga"
$ws.Range("H495").Value = "Sin especificar"
$ws.Range("I495").Value = "Primera"
$ws.Range("J495").Value = 1000
$ws.Range("K495").Value = 1000
$ws.Range("L495").Value = 1100
$ws.Range("M495").Value = 1050
$ws.Range("N495").Value = "`$/paquete 5 unidades"
$ws.Range("O495").Value = "Región Metropolitana"
$ws.Range("P495").Value = 210
$ws.Range("Q495").Value = 5
$ws.Range("R495").Value = "Hortaliza"
